$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.489076256752014
$ws.Range("B1").Value = 2.819122314453125
$ws.Range("C1").Value = 6.931576728820801
$ws.Range("D1").Value = 1.745533466339111
$ws.Range("E1").Value = 0.8950828313827515
